$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '70.285.99'
    'E2' = '  -2.83%  '
    'D3' = '2.521.15'
    'E3' = '  -5.15%  '
    'E4' = '  -0.06%  '
    'D5' = '575.09'
    'E5' = '  -3.77%  '
    'D6' = '169.76'
    'E6' = '  -2.92%  '
    'E7' = '  +0.08%  '
    'E8' = '  -2.79%  '
    'D9' = '2.520.03'
    'E9' = '  -5.18%  '
    'E10' = '  -6.02%  '
    'E11' = '  -0.65%  '
    'D12' = '0.344'
    'E12' = '  -4.09%  '
    'D13' = '4.83'
    'E13' = '  -3.45%  '
    'D14' = '2.983.73'
    'E14' = '  -5.20%  '
    'D15' = '70.124.06'
    'E15' = '  -3.04%  '
    'E16' = '  -3.41%  '
    'D17' = '25.03'
    'E17' = '  -4.96%  '
    'D18' = '2.516.34'
    'E18' = '  -5.24%  '
    'D19' = '11.53'
    'E19' = '  -6.02%  '
    'D20' = '7.59'
    'E20' = '  -8.23%  '
    'D21' = '355.14'
    'E21' = '  -4.07%  '
    'D22' = '3.95'
    'E22' = '  -5.49%  '
    'D23' = '1.99'
    'E23' = '  -3.14%  '
    'E24' = '  +0.05%  '
    'D25' = '69.03'
    'E25' = '  -4.33%  '
    'E26' = '  -5.54%  '
    'D27' = '9.21'
    'E27' = '  -5.69%  '
    'D28' = '2.649.46'
    'E28' = '  -5.41%  '
    'D29' = '0.998'
    'E29' = '  -0.11%  '
    'D30' = '0.0₃0911'
    'E30' = '  -6.60%  '
    'D31' = '7.85'
    'E31' = '  -3.42%  '
    'D32' = '483.13'
    'E32' = '  -3.54%  '
    'D33' = '1.30'
    'E33' = '  +0.12%  '
    'E34' = '  -3.59%  '
    'D35' = '0.999'
    'E35' = '  -0.10%  '
    'E36' = '  +4.18%  '
    'D37' = '156.03'
    'E37' = '  -4.19%  '
    'E38' = '  -0.34%  '
    'D39' = '18.58'
    'E39' = '  -4.86%  '
    'E40' = '  +0.05%  '
    'B41' = 'PolygonEcosystemToken'
    'C41' = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
    'D41' = '0.321'
    'E41' = '  -3.45%  '
    'D42' = '1.65'
    'E42' = '  -7.36%  '
    'B43' = 'RenderToken'
    'C43' = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
    'D43' = '4.74'
    'E43' = '  -5.15%  '
    'B44' = 'ImmutableX'
    'C44' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D44' = '1.23'
    'E44' = '  -11.33%  '
    'D45' = '2.39'
    'E45' = '  -8.22%  '
    'D46' = '38.31'
    'E46' = '  -3.03%  '
    'D47' = '142.40'
    'E47' = '  -8.90%  '
    'E48' = '  -5.90%  '
    'D49' = '0.527'
    'E49' = '  -5.69%  '
    'D50' = '1.62'
    'E50' = '  -6.34%  '
    'D51' = '0.599'
    'E51' = '  -0.92%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}

Write-Output "Applied $($updates.Count) cell updates"